$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ram"

$ws.Rows.Item(4).Delete()
